$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.132.99'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.655.26'
$ws.Range('E3').Value = '  -0.76%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.23'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5287'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2609'
$ws.Range('E8').Value = '  -2.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06340'
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.45'
$ws.Range('E10').Value = '  -2.66%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.496'
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('D13').Value = '1.653.64'
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5475'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '0.0₅8170'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.43'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = '26.133.99'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.560'
$ws.Range('E19').Value = '  -2.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.14'
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('E21').Value = '  -0.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.028'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '141.84'
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1249'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.24'
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.437'
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05945'
$ws.Range('E29').Value = '  -4.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.282'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.517'
$ws.Range('E31').Value = '  -2.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.246'
$ws.Range('E32').Value = '  -1.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.572'
$ws.Range('E33').Value = '  -3.72%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9502'
$ws.Range('E34').Value = '  -2.47%  '
$ws.Range('B35').Value = 'MXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.790'
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.409'
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5661'
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01611'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.809'
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8485'
$ws.Range('E40').Value = '  -1.13%  '
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '102.77'
$ws.Range('E42').Value = '  +2.35%  '
$ws.Range('D43').Value = '1.023.93'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.20'
$ws.Range('E45').Value = '  -0.95%  '
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4289'
$ws.Range('E47').Value = '  +1.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.478'
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05152'
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.819'
$ws.Range('E50').Value = '  -3.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.09705'
$ws.Range('E51').Value = '  -0.97%  '
